$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.305.40"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "1.690.03"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'217.90"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "'0.5372"
$ws.Range("E6").Value = "  +2.32%  "
$ws.Range("D7").Value = "'1.008"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.2729"
$ws.Range("E8").Value = "  +1.31%  "
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("D10").Value = "'21.56"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").Value = "'0.07668"
$ws.Range("E11").Value = "  +1.87%  "
$ws.Range("D12").Value = "1.696.86"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "'4.528"
$ws.Range("E13").Value = "  -0.04%  "
$ws.Range("D14").Value = "'0.5780"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "'0.000008372"
$ws.Range("E15").Value = "  -1.68%  "
$ws.Range("E16").Value = "  +3.37%  "
$ws.Range("D17").Value = "26.359.36"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "'4.900"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D21").Value = "'190.48"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").Value = "'6.254"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'149.05"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("D26").Value = "'7.847"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "'15.88"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").Value = "'0.06243"
$ws.Range("E28").Value = "  -3.37%  "
$ws.Range("D29").Value = "'1.370"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D31").Value = "'3.594"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "'3.580"
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("D33").Value = "'1.670"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").Value = "'1.030"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "'0.6144"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").Value = "'2.420"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").Value = "'2.761"
$ws.Range("E37").Value = "  +1.02%  "
$ws.Range("D38").Value = "'0.01653"
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("D39").Value = "1.108.13"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "'6.112"
$ws.Range("E40").Value = "  -5.09%  "
$ws.Range("D41").Value = "'0.8803"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").Value = "'101.34"
$ws.Range("D44").Value = "1.841.75"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").Value = "'0.00000000111"
$ws.Range("E45").Value = "  -1.51%  "
$ws.Range("D46").Value = "'57.59"
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").Value = "'8.122"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("D48").Value = "'1.001"
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("D50").Value = "'0.4299"
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("D51").Value = "'6.030"
$ws.Range("E51").Value = "  -0.86%  "
